$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'90.842.18"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.12%  '

$ws.Range('D3').Value = "'3.155.02"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +2.83%  '

$ws.Range('E4').Value = '  -0.47%  '

$ws.Range('D5').Value = "'215.84"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.68%  '

$ws.Range('D6').Value = "'627.84"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.91%  '

$ws.Range('E7').Value = '  +31.96%  '

$ws.Range('D8').Value = "'0.369"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.00%  '

$ws.Range('D9').Value = "'0.999"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.00%  '

$ws.Range('D10').Value = "'3.154.10"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +2.76%  '

$ws.Range('D11').Value = "'0.757"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +11.17%  '

$ws.Range('E12').Value = '  +7.34%  '

$ws.Range('D13').Value = "'5.75"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +7.13%  '

$ws.Range('D14').Value = "'0.0000246"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.28%  '

$ws.Range('D15').Value = "'34.92"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +6.57%  '

$ws.Range('D16').Value = "'90.693.15"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.06%  '

$ws.Range('D17').Value = "'3.739.00"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.98%  '

$ws.Range('D18').Value = "'3.135.84"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.27%  '

$ws.Range('D19').Value = "'3.75"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +8.87%  '

$ws.Range('D20').Value = "'14.63"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +6.84%  '

$ws.Range('D21').Value = "'476.05"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +10.59%  '

$ws.Range('D22').Value = "'0.0000211"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -4.96%  '

$ws.Range('D23').Value = "'9.18"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +9.22%  '

$ws.Range('D24').Value = "'5.27"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +3.83%  '

$ws.Range('D25').Value = "'96.00"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +15.17%  '

$ws.Range('D26').Value = "'5.94"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +8.63%  '

$ws.Range('D27').Value = "'12.35"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +4.63%  '

$ws.Range('D28').Value = "'3.329.85"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +3.86%  '

$ws.Range('E29').Value = '  -0.08%  '

$ws.Range('E30').Value = '  -1.82%  '

$ws.Range('D31').Value = "'9.28"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +7.30%  '

$ws.Range('E32').Value = '  -6.03%  '

$ws.Range('D33').Value = "'27.63"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +20.31%  '

$ws.Range('D34').Value = "'0.205"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +48.91%  '

$ws.Range('D35').Value = "'518.86"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.19%  '

$ws.Range('E36').Value = '  +6.02%  '

$ws.Range('D37').Value = "'0.145"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +7.00%  '

$ws.Range('B38').Value = 'RenderToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D38').Value = "'6.94"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.63%  '

$ws.Range('B39').Value = 'dogwifhat'
$ws.Range('C39').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D39').Value = "'3.61"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -4.90%  '

$ws.Range('D40').Value = "'1.31"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +3.70%  '

$ws.Range('D41').Value = "'0.0924"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +28.50%  '

$ws.Range('D42').Value = "'22.22"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.28%  '

$ws.Range('B43').Value = 'PolygonEcosystemToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D43').Value = "'0.422"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +15.58%  '

$ws.Range('B44').Value = 'FirstDigitalUSD'
$ws.Range('C44').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D44').Value = "'1.00"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.05%  '

$ws.Range('D45').Value = "'1.98"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +6.60%  '

$ws.Range('E46').Value = '  +0.03%  '

$ws.Range('D47').Value = "'0.733"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +21.93%  '

$ws.Range('D48').Value = "'151.16"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +5.13%  '

$ws.Range('D49').Value = "'4.69"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +10.83%  '

$ws.Range('B50').Value = 'ImmutableX'
$ws.Range('C50').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D50').Value = "'1.37"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +10.68%  '

$ws.Range('B51').Value = 'OKB'
$ws.Range('C51').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D51').Value = "'45.63"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +4.53%  '
